$d = $word.ActiveDocument

# Locate the paragraph that ends with "...as List.cshtml."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Created Item Controller\. Created Views for Item as List\.cshtml\.") {
        $target = $p
        break
    }
}

# Insert a new paragraph right after the found paragraph.
$newPara = $target.Range.InsertParagraphAfter()

# The newly created paragraph is the one right after $target.
$newP = $target.Next()
$newP.Range.Text = "Created ViewModels to display data instead of View Bag."
